# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted before
# the existing "Late" column (column N), shifting "Late" and "Outstanding"
# one column to the right. The "Repayment Schedule" tab also becomes the
# active/selected sheet (it was "Transactions" before), with cell P12
# selected on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N (pushes "Late"/"Outstanding" right).
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with P12 selected.
$ws.Activate() | Out-Null
$ws.Range("P12").Select() | Out-Null
